# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per latest data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 111112120
$ws.Range("I32").Value = 500000480
$ws.Range("J32").Value = 1160.7142
$ws.Range("K32").Value = 500000480
$ws.Range("L32").Value = 1160.7142
$ws.Range("M32").Value = -500000154
$ws.Range("N32").Value = -1812.7142

$ws.Range("H106").Value = 3272.5715
$ws.Range("I106").Value = 3272.5715
$ws.Range("K106").Value = 3272.5715
$ws.Range("M106").Value = -2641.5715

$ws.Range("H125").Value = 3434.6667
$ws.Range("I125").Value = 6016
$ws.Range("J125").Value = 2697.1428
$ws.Range("K125").Value = 54144
$ws.Range("L125").Value = 24274.2852
$ws.Range("M125").Value = -51684
$ws.Range("N125").Value = -29194.2852

$ws.Range("H129").Value = 878.075
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 887.7692
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 2663.3076
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -12663.3076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 41713228
$ws.Range("I110").Value = 43526804
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 43526804
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = -43524759
$ws.Range("N110").Value = -5090

$ws.Range("H132").Value = 2984.9792
$ws.Range("I132").Value = 2923.6428
$ws.Range("J132").Value = 3414.3333
$ws.Range("K132").Value = 8770.928400000001
$ws.Range("L132").Value = 10242.9999
$ws.Range("M132").Value = -6240.928400000001
$ws.Range("N132").Value = -15302.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H64").Value = 412.6154
$ws.Range("J64").Value = 511
$ws.Range("L64").Value = 511
$ws.Range("N64").Value = -961

$ws.Range("H67").Value = 412.6154
$ws.Range("J67").Value = 511
$ws.Range("L67").Value = 511
$ws.Range("N67").Value = -2071

$ws.Range("H96").Value = 15950
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492

$ws.Range("H107").Value = 38462224
$ws.Range("I107").Value = 71429190
$ws.Range("J107").Value = 765.4167
$ws.Range("K107").Value = 71429190
$ws.Range("L107").Value = 765.4167
$ws.Range("M107").Value = -71427270
$ws.Range("N107").Value = -4605.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 9636.857
$ws.Range("I122").Value = 9889.5
$ws.Range("J122").Value = 9300
$ws.Range("K122").Value = 29668.5
$ws.Range("L122").Value = 27900
$ws.Range("M122").Value = -27218.5
$ws.Range("N122").Value = -32800

$ws.Range("H132").Value = 4599.3335
$ws.Range("I132").Value = 4477.1113
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 13431.3339
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -10901.3339
$ws.Range("N132").Value = -21057.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 12415.944
$ws.Range("I5").Value = 1555.6
$ws.Range("J5").Value = 16593
$ws.Range("K5").Value = 4666.799999999999
$ws.Range("L5").Value = 49779
$ws.Range("M5").Value = -4554.799999999999
$ws.Range("N5").Value = -50003

$ws.Range("H8").Value = 295.83334
$ws.Range("I8").Value = 295.83334
$ws.Range("K8").Value = 887.5000200000001
$ws.Range("M8").Value = -748.5000200000001

$ws.Range("H107").Value = 2167.9
$ws.Range("J107").Value = 2222.375
$ws.Range("L107").Value = 6667.125
$ws.Range("N107").Value = -10507.125

$ws.Range("H116").Value = 2155.3333
$ws.Range("J116").Value = 3750
$ws.Range("L116").Value = 11250
$ws.Range("N116").Value = -18134

$ws.Range("H117").Value = 2256.4285
$ws.Range("J117").Value = 2952.8
$ws.Range("L117").Value = 8858.400000000001
$ws.Range("N117").Value = -15742.4

$ws.Range("H121").Value = 10937.286
$ws.Range("I121").Value = 6939.7144
$ws.Range("J121").Value = 14934.857
$ws.Range("K121").Value = 20819.1432
$ws.Range("L121").Value = 44804.571
$ws.Range("M121").Value = -19509.1432
$ws.Range("N121").Value = -47424.571

$ws.Range("H122").Value = 4243.2593
$ws.Range("I122").Value = 322.1579
$ws.Range("K122").Value = 2899.4211
$ws.Range("M122").Value = -449.4211

$ws.Range("H131").Value = 785.0700000000001
$ws.Range("I131").Value = 390.69232
$ws.Range("J131").Value = 844
$ws.Range("K131").Value = 1172.07696
$ws.Range("L131").Value = 2532
$ws.Range("M131").Value = 3867.92304
$ws.Range("N131").Value = -12612

$ws.Range("H135").Value = 12415.944
$ws.Range("I135").Value = 1555.6
$ws.Range("J135").Value = 16593
$ws.Range("K135").Value = 14000.4
$ws.Range("L135").Value = 149337
$ws.Range("M135").Value = -11465.4
$ws.Range("N135").Value = -154407

$ws.Range("H138").Value = 2313.25
$ws.Range("I138").Value = 2022.7142
$ws.Range("J138").Value = 2720
$ws.Range("K138").Value = 6068.142599999999
$ws.Range("L138").Value = 8160
$ws.Range("M138").Value = -928.1425999999992
$ws.Range("N138").Value = -18440

$ws.Range("H140").Value = 2358.3333
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 2358.3333
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").Value = 7074.999899999999
$ws.Range("N140").Value = -17434.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 297.25
$ws.Range("I2").Value = 281.5
$ws.Range("J2").Value = 376
$ws.Range("K2").Value = 281.5
$ws.Range("L2").Value = 376
$ws.Range("M2").Value = -168.5
$ws.Range("N2").Value = -602

$ws.Range("H113").Value = 2321.4666
$ws.Range("I113").Value = 2744.5715
$ws.Range("K113").Value = 2744.5715
$ws.Range("M113").Value = -574.5715

$ws.Range("H132").Value = 2493.8965
$ws.Range("I132").Value = 1878.3889
$ws.Range("J132").Value = 3501.0908
$ws.Range("K132").Value = 5635.1667
$ws.Range("L132").Value = 10503.2724
$ws.Range("M132").Value = -3105.1667
$ws.Range("N132").Value = -15563.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 78957.84
$ws.Range("I40").Value = 334401.34
$ws.Range("J40").Value = 2324.8
$ws.Range("K40").Value = 334401.34
$ws.Range("L40").Value = 2324.8
$ws.Range("M40").Value = -334265.34
$ws.Range("N40").Value = -2596.8

$ws.Range("H46").Value = 3398
$ws.Range("I46").Value = 540
$ws.Range("J46").Value = 4214.5713
$ws.Range("K46").Value = 540
$ws.Range("L46").Value = 4214.5713
$ws.Range("M46").Value = -352
$ws.Range("N46").Value = -4590.5713

$ws.Range("H61").Value = 2785.4443
$ws.Range("I61").Value = 1917.8
$ws.Range("K61").Value = 1917.8
$ws.Range("M61").Value = -1715.8

$ws.Range("H96").Value = 15298.8
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H113").Value = 2785.4443
$ws.Range("I113").Value = 1917.8
$ws.Range("K113").Value = 1917.8
$ws.Range("M113").Value = 252.2

$ws.Range("H122").Value = 2535.7693
$ws.Range("I122").Value = 1825.875
$ws.Range("J122").Value = 3671.6
$ws.Range("K122").Value = 5477.625
$ws.Range("L122").Value = 11014.8
$ws.Range("M122").Value = -3027.625
$ws.Range("N122").Value = -15914.8

$ws.Range("H132").Value = 3634.7083
$ws.Range("I132").Value = 3834.7222
$ws.Range("J132").Value = 3034.6667
$ws.Range("K132").Value = 11504.1666
$ws.Range("L132").Value = 9104.000100000001
$ws.Range("M132").Value = -8974.1666
$ws.Range("N132").Value = -14164.0001
